$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Test 5" results column is inserted between the existing "Test 4"
# column (G) and the "Average" column, which slides from H to I.
# Row 1 (the merged title band) has no content in this area, so it is left
# untouched - only the header row and the two data rows actually move.

# Header row: H2 becomes "Test 5", I2 becomes the (moved) "Average" label.
$ws.Range("H2").Value = "Test 5"
$ws.Range("I2").Value = "Average"

# Data row 3: the old average formula (over D3:G3) moves from H3 to I3;
# H3 is left blank since there's no "Test 5" value for this row yet.
$ws.Range("I3").Formula = "=AVERAGE(D3:G3)"
$ws.Range("H3").Clear()

# Data row 4: H4 now holds the literal "Test 5" result, and I4 carries the
# updated average formula spanning the new column D4:H4.
$ws.Range("H4").Value = 76.08
$ws.Range("I4").Formula = "=AVERAGE(D4:H4)"

# Match the author's final selection/scroll position.
$ws.Range("I6").Select()
